$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 878.7778
$ws.Range("I43").Value = 599
$ws.Range("J43").Value = 913.75
$ws.Range("K43").Value = 599
$ws.Range("L43").Value = 913.75
$ws.Range("M43").Value = -530
$ws.Range("N43").Value = -1051.75

$ws.Range("H100").Value = 3226.6667
$ws.Range("J100").Value = 3900
$ws.Range("L100").Value = 3900
$ws.Range("N100").Value = -4982

$ws.Range("H103").Value = 111111384
$ws.Range("J103").Value = 420
$ws.Range("L103").Value = 1260
$ws.Range("N103").Value = -2432

$ws.Range("H129").Value = 176425.94
$ws.Range("I129").Value = 350
$ws.Range("J129").Value = 182828.69
$ws.Range("K129").Value = 1050
$ws.Range("L129").Value = 548486.0700000001
$ws.Range("M129").Value = 3950
$ws.Range("N129").Value = -558486.0700000001

$ws.Range("H132").Value = 2674.25
$ws.Range("I132").Value = 2781
$ws.Range("K132").Value = 8343
$ws.Range("M132").Value = -5813

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7303.1284
$ws.Range("I32").Value = 6023.393
$ws.Range("J32").Value = 12422.071
$ws.Range("K32").Value = 6023.393
$ws.Range("L32").Value = 12422.071
$ws.Range("M32").Value = -5736.393
$ws.Range("N32").Value = -12996.071

$ws.Range("H132").Value = 10352.807
$ws.Range("I132").Value = 1334.24
$ws.Range("K132").Value = 4002.72
$ws.Range("M132").Value = -1472.72

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 714.4074000000001
$ws.Range("J94").Value = 786.6667
$ws.Range("L94").Value = 786.6667
$ws.Range("N94").Value = -1688.6667

$ws.Range("H134").Value = 3510.2954
$ws.Range("I134").Value = 3464.946
$ws.Range("K134").Value = 10394.838
$ws.Range("M134").Value = -7859.838

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5000
$ws.Range("J4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5224

$ws.Range("H31").Value = 2924.7454
$ws.Range("J31").Value = 4562.76
$ws.Range("L31").Value = 4562.76
$ws.Range("N31").Value = -5152.76

$ws.Range("H34").Value = 2924.7454
$ws.Range("J34").Value = 4562.76
$ws.Range("L34").Value = 4562.76
$ws.Range("N34").Value = -4966.76

$ws.Range("H58").Value = 31365.176
$ws.Range("I58").Value = 2098
$ws.Range("K58").Value = 2098
$ws.Range("M58").Value = -1895

$ws.Range("H132").Value = 1548.9166
$ws.Range("I132").Value = 1119.2885
$ws.Range("J132").Value = 4341.5
$ws.Range("K132").Value = 3357.8655
$ws.Range("L132").Value = 13024.5
$ws.Range("M132").Value = -827.8655000000003
$ws.Range("N132").Value = -18084.5

$ws.Range("H134").Value = 970.2857
$ws.Range("I134").Value = 891.0769
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2673.2307
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -138.2307000000001
$ws.Range("N134").Value = -11070

$ws.Range("H136").Value = 31365.176
$ws.Range("I136").Value = 2098
$ws.Range("K136").Value = 6294
$ws.Range("M136").Value = -3744

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 922.6053000000001
$ws.Range("I5").Value = 728.0741
$ws.Range("J5").Value = 1400.091
$ws.Range("K5").Value = 2184.2223
$ws.Range("L5").Value = 4200.272999999999
$ws.Range("M5").Value = -2072.2223
$ws.Range("N5").Value = -4424.272999999999

$ws.Range("H131").Value = 681.09
$ws.Range("J131").Value = 707.7912
$ws.Range("L131").Value = 2123.3736
$ws.Range("N131").Value = -12203.3736

$ws.Range("H134").Value = 3372.3684
$ws.Range("I134").Value = 2051.8
$ws.Range("J134").Value = 4839.6665
$ws.Range("K134").Value = 6155.400000000001
$ws.Range("L134").Value = 14518.9995
$ws.Range("M134").Value = -1085.400000000001
$ws.Range("N134").Value = -24658.9995

$ws.Range("H135").Value = 922.6053000000001
$ws.Range("I135").Value = 728.0741
$ws.Range("J135").Value = 1400.091
$ws.Range("K135").Value = 6552.6669
$ws.Range("L135").Value = 12600.819
$ws.Range("M135").Value = -4017.6669
$ws.Range("N135").Value = -17670.819

$ws.Range("H139").Value = 2281
$ws.Range("I139").Value = 1471.3043
$ws.Range("J139").Value = 3444.9375
$ws.Range("K139").Value = 4413.9129
$ws.Range("L139").Value = 10334.8125
$ws.Range("M139").Value = 726.0870999999997
$ws.Range("N139").Value = -20614.8125

$ws.Range("H140").Value = 1990.5312
$ws.Range("I140").Value = 957.0526
$ws.Range("J140").Value = 3501
$ws.Range("K140").Value = 2871.1578
$ws.Range("L140").Value = 10503
$ws.Range("M140").Value = 2308.8422
$ws.Range("N140").Value = -20863

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 7337.0557
$ws.Range("I113").Value = 8459.071
$ws.Range("K113").Value = 8459.071
$ws.Range("M113").Value = -6289.071

$ws.Range("H126").Value = 2984.1277
$ws.Range("I126").Value = 2381.353
$ws.Range("J126").Value = 4560.615
$ws.Range("K126").Value = 7144.059
$ws.Range("L126").Value = 13681.845
$ws.Range("M126").Value = -4674.059
$ws.Range("N126").Value = -18621.845

$ws.Range("H132").Value = 18986.719
$ws.Range("J132").Value = 74379.57000000001
$ws.Range("L132").Value = 223138.71
$ws.Range("N132").Value = -228198.71

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2703.125
$ws.Range("I22").Value = 3128.5386
$ws.Range("J22").Value = 859.6667
$ws.Range("K22").Value = 3128.5386
$ws.Range("L22").Value = 859.6667
$ws.Range("M22").Value = -2833.5386
$ws.Range("N22").Value = -1449.6667

$ws.Range("H27").Value = 2703.125
$ws.Range("I27").Value = 3128.5386
$ws.Range("J27").Value = 859.6667
$ws.Range("K27").Value = 3128.5386
$ws.Range("L27").Value = 859.6667
$ws.Range("M27").Value = -3021.5386
$ws.Range("N27").Value = -1073.6667

$ws.Range("H82").Value = 1451.3334
$ws.Range("I82").Value = 1270.3334
$ws.Range("K82").Value = 1270.3334
$ws.Range("M82").Value = -909.3334

$ws.Range("H85").Value = 1451.3334
$ws.Range("I85").Value = 1270.3334
$ws.Range("K85").Value = 1270.3334
$ws.Range("M85").Value = -22.33339999999998

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws.Range("H136").Value = 1711.5238
$ws.Range("I136").Value = 1497.1
$ws.Range("K136").Value = 4491.299999999999
$ws.Range("M136").Value = -1941.299999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 4000
$ws.Range("J15").Value = 4000
$ws.Range("L15").Value = 4000
$ws.Range("N15").Value = -4576

$ws.Range("H58").Value = 16000
$ws.Range("I58").Value = 5000
$ws.Range("J58").Value = 17833.334
$ws.Range("K58").Value = 5000
$ws.Range("L58").Value = 17833.334
$ws.Range("M58").Value = -4692
$ws.Range("N58").Value = -18449.334

$ws.Range("H62").Value = 4643.2856
$ws.Range("J62").Value = 4643.2856
$ws.Range("L62").Value = 4643.2856
$ws.Range("N62").Value = -5891.2856

$ws.Range("H65").Value = 4643.2856
$ws.Range("J65").Value = 4643.2856
$ws.Range("L65").Value = 23216.428
$ws.Range("N65").Value = -29456.428
